$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list values in column D
$ws.Range("D32").Value = 10881.07
$ws.Range("D33").Value = 8642.678
$ws.Range("D34").Value = 8621.951999999999
$ws.Range("D35").Value = 12404.421
